$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second header row (it gets merged into row 1 as new column headers)
$ws.Rows(2).Delete()

# Clear any leftover formatting on the left-hand header cells so they come
# back to the default (unstyled) look used by the new headers.
$ws.Range("A1:E1").Style = "Normal"

# Rewrite the header row (row 1) with the new column layout.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The right-hand headers (F1:K1) carry the 9pt Arial font style used
# elsewhere in the sheet, applied via a transient named style so a fresh
# cell-format record (font-only, no explicit number format) gets created.
$tempStyle = $wb.Styles.Add("TempHeaderStyle")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TempHeaderStyle"
$wb.Styles("TempHeaderStyle").Delete()

$ws.Range("A2:K2").Select() | Out-Null
